# Fix page numbers in presentations
#
# Every content slide (2..28) carries a small "page number" text box
# (originally showing the stale "2/20") in the bottom-right corner of the
# slide. This script updates each one to show "<slide index>/28" (the
# deck has 28 slides total), and - for the slides whose new label is
# wider than the old "2/20" - grows/repositions the box the same way the
# author's edit did (new Left/Width, and the box brought in front of the
# other shapes on the slide).

$p = $ppt.ActivePresentation
$EMU_PER_PT = 12700

# slide index -> list of run texts that, concatenated, form "<n>/28"
$runsMap = @{}
$runsMap[2]  = @("2/28")
$runsMap[3]  = @("3", "/28")
$runsMap[4]  = @("4", "/28")
$runsMap[5]  = @("5", "/28")
$runsMap[6]  = @("6", "/28")
$runsMap[7]  = @("7", "/28")
$runsMap[8]  = @("8", "/28")
$runsMap[9]  = @("9", "/28")
$runsMap[10] = @("10", "/28")
$runsMap[11] = @("11", "/28")
$runsMap[12] = @("12", "/28")
$runsMap[13] = @("13", "/28")
$runsMap[14] = @("14", "/28")
$runsMap[15] = @("15", "/28")
$runsMap[16] = @("16", "/28")
$runsMap[17] = @("17", "/28")
$runsMap[18] = @("18", "/28")
$runsMap[19] = @("19", "/28")
$runsMap[20] = @("2", "0", "/28")
$runsMap[21] = @("21/28")
$runsMap[22] = @("22/28")
$runsMap[23] = @("23", "/28")
$runsMap[24] = @("24/28")
$runsMap[25] = @("25/28")
$runsMap[26] = @("26/28")
$runsMap[27] = @("27/28")
$runsMap[28] = @("28/28")

# slide index -> new box width in EMU, for the slides whose label grew
# wider than the original 4-character "2/20" (i.e. every two-digit page
# number, 10..28). Everything else keeps its original box geometry.
$widthMap = @{}
$widthMap[10] = 761747
$widthMap[11] = 744627
$widthMap[12] = 761747
$widthMap[13] = 761747
$widthMap[14] = 761747
$widthMap[15] = 761747
$widthMap[16] = 761747
$widthMap[17] = 761747
$widthMap[18] = 761747
$widthMap[19] = 761747
$widthMap[20] = 761747
$widthMap[21] = 761747
$widthMap[22] = 761747
$widthMap[23] = 761747
$widthMap[24] = 761747
$widthMap[25] = 761747
$widthMap[26] = 761747
$widthMap[27] = 761747
$widthMap[28] = 761747

$newLeftEmu = 8346757

for ($slideIdx = 2; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    # Locate the stale page-number text box by its exact text ("2/20"),
    # rather than assuming a fixed shape index (it varies slide to slide).
    $target = $null
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "2/20") {
                $target = $shp
            }
        }
    }

    if ($target -eq $null) {
        continue
    }

    $runs = $runsMap[$slideIdx]
    if ($runs -eq $null) {
        continue
    }

    $tr = $target.TextFrame.TextRange
    $tr.Text = $runs[0]
    for ($r = 1; $r -lt $runs.Count; $r++) {
        $unused = $tr.InsertAfter($runs[$r])
    }

    if ($widthMap.ContainsKey($slideIdx)) {
        $target.Left = $newLeftEmu / $EMU_PER_PT
        $target.Width = $widthMap[$slideIdx] / $EMU_PER_PT
    }

    # The slides that grew a wider label (everything except the very
    # first two-digit one, slide 10) also had their box brought to the
    # front of the z-order in the authored edit.
    if ($widthMap.ContainsKey($slideIdx) -and $slideIdx -ne 10) {
        $target.ZOrder(0)
    }
}
